$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Quarterly unemployment-rate table refresh ------------------------------
# One new quarter (01/01/2024) was appended to each of the three region
# blocks (Brasil: rows 2-22, Nordeste: rows 23-43, Sergipe: rows 44-64).
# Every "Trimestre" (column C) and "Valor" (column D) below it shifts up by
# one row to make room, and the newest quarter's value lands on the last
# row of each block.
#
# Each table row below is: (row number, new Trimestre text, new Valor or $null)
$updates = @(
    @(2, "01/01/2019", 12.8),
    @(3, "01/04/2019", 12.1),
    @(4, "01/07/2019", 11.9),
    @(5, "01/10/2019", 11.1),
    @(6, "01/01/2020", 12.4),
    @(7, "01/04/2020", 13.6),
    @(8, "01/07/2020", 14.9),
    @(9, "01/10/2020", 14.2),
    @(10, "01/01/2021", 14.9),
    @(11, "01/04/2021", 14.2),
    @(12, "01/07/2021", 12.6),
    @(13, "01/10/2021", 11.1),
    @(14, "01/01/2022", 11.1),
    @(15, "01/04/2022", 9.300000000000001),
    @(16, "01/07/2022", 8.699999999999999),
    @(17, "01/10/2022", 7.9),
    @(18, "01/01/2023", 8.800000000000001),
    @(19, "01/04/2023", 8),
    @(20, "01/07/2023", 7.7),
    @(21, "01/10/2023", 7.4),
    @(22, "01/01/2024", 7.9),
    @(23, "01/01/2019", 15.4),
    @(24, "01/04/2019", 14.8),
    @(25, "01/07/2019", 14.6),
    @(26, "01/10/2019", 13.8),
    @(27, "01/01/2020", 15.8),
    @(28, "01/04/2020", $null),
    @(29, "01/07/2020", $null),
    @(30, "01/10/2020", $null),
    @(31, "01/01/2021", $null),
    @(32, "01/04/2021", $null),
    @(33, "01/07/2021", $null),
    @(34, "01/10/2021", $null),
    @(35, "01/01/2022", $null),
    @(36, "01/04/2022", 12.7),
    @(37, "01/07/2022", 12),
    @(38, "01/10/2022", 10.9),
    @(39, "01/01/2023", 12.2),
    @(40, "01/04/2023", 11.3),
    @(41, "01/07/2023", 10.8),
    @(42, "01/10/2023", 10.4),
    @(43, "01/01/2024", 11.1),
    @(44, "01/01/2019", 15.5),
    @(45, "01/04/2019", 15.4),
    @(46, "01/07/2019", 14.8),
    @(47, "01/10/2019", 15),
    @(48, "01/01/2020", 15.8),
    @(49, "01/04/2020", $null),
    @(50, "01/07/2020", $null),
    @(51, "01/10/2020", $null),
    @(52, "01/01/2021", $null),
    @(53, "01/04/2021", $null),
    @(54, "01/07/2021", $null),
    @(55, "01/10/2021", $null),
    @(56, "01/01/2022", $null),
    @(57, "01/04/2022", 12.7),
    @(58, "01/07/2022", 12.1),
    @(59, "01/10/2022", 11.9),
    @(60, "01/01/2023", 11.9),
    @(61, "01/04/2023", 10.3),
    @(62, "01/07/2023", 9.800000000000001),
    @(63, "01/10/2023", 11.2),
    @(64, "01/01/2024", 10)
)

# Column C holds quarter labels like "01/01/2019" which look like dates.
# Force the column to Text first so Excel stores them as the literal
# strings from the source data instead of silently coercing them to date
# serial numbers.
$ws.Range("C2:C64").NumberFormat = "@"

foreach ($row in $updates) {
    $r = $row[0]
    $trimestre = $row[1]
    $valor = $row[2]

    $ws.Cells.Item($r, 3).Value = $trimestre

    if ($null -eq $valor) {
        $ws.Cells.Item($r, 4).ClearContents()
    } else {
        $ws.Cells.Item($r, 4).Value = $valor
    }
}
